# excel_writer: also include totals for the balance columns
#
# The "Gesamtergebnis" (overall result) sheet has a Total row (row 3) that
# previously left the two balance columns (Startguthaben/Endsaldo -> C/D)
# as "N/A" placeholders. They should now carry a numeric total (0, matching
# the other totals in that row) just like every other summed column.
#
# Alongside that data fix, the workbook's remembered selection/active-sheet
# state (as last saved by the tool) moves: the day-results sheet is no
# longer the selected tab, the month-results sheet's selection collapses to
# a single cell, and the overall-result sheet becomes the active tab with
# its selection on C5.

$wb = $excel.ActiveWorkbook

$wsDay   = $wb.Worksheets.Item("Tagesergebnisse")
$wsMonth = $wb.Worksheets.Item("Monatsergebnisse")
$wsTotal = $wb.Worksheets.Item("Gesamtergebnis")

# --- Data fix: totals for the balance columns on the overall-result sheet ---
$wsTotal.Range("C3").Value = 0
$wsTotal.Range("D3").Value = 0

# Row 3 (the Total row) recalculates to the shorter "numeric row" height
# used elsewhere in the workbook once it no longer holds the taller N/A text.
$wsTotal.Rows.Item(3).RowHeight = 13.8

# --- Selection / active sheet bookkeeping ---

# Day-results sheet keeps its A2 selection but is no longer the active tab.
[void]$wsDay.Range("A2").Select()

# Month-results sheet's selection collapses down to the single cell A6.
[void]$wsMonth.Range("A6").Select()

# Overall-result sheet becomes selected/active, with the cursor on C5.
[void]$wsTotal.Range("C5").Select()
[void]$wsTotal.Activate()
